$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "Prevention"
